# Update division problems in the three-digit_number_divided_by_one-digit_number worksheet
$d = $word.ActiveDocument

$d.Content.Find.Execute("995÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "725÷7=", 2)
$d.Content.Find.Execute("692÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "650÷6=", 2)
$d.Content.Find.Execute("438÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "928÷4=", 2)
$d.Content.Find.Execute("998÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "294÷4=", 2)
$d.Content.Find.Execute("576÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "644÷7=", 2)
$d.Content.Find.Execute("267÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "784÷8=", 2)
$d.Content.Find.Execute("109÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "979÷9=", 2)
$d.Content.Find.Execute("108÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "763÷7=", 2)
$d.Content.Find.Execute("406÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "732÷2=", 2)
$d.Content.Find.Execute("567÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "111÷4=", 2)
$d.Content.Find.Execute("696÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "487÷7=", 2)
$d.Content.Find.Execute("557÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "190÷7=", 2)
$d.Content.Find.Execute("922÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "322÷7=", 2)
$d.Content.Find.Execute("928÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "856÷4=", 2)
$d.Content.Find.Execute("842÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "801÷3=", 2)
$d.Content.Find.Execute("715÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "171÷5=", 2)
$d.Content.Find.Execute("325÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "923÷7=", 2)
$d.Content.Find.Execute("129÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "242÷2=", 2)
$d.Content.Find.Execute("288÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "231÷8=", 2)
$d.Content.Find.Execute("924÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "765÷7=", 2)
$d.Content.Find.Execute("483÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "537÷3=", 2)
$d.Content.Find.Execute("511÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "867÷7=", 2)
$d.Content.Find.Execute("415÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "342÷8=", 2)
$d.Content.Find.Execute("447÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "847÷8=", 2)
$d.Content.Find.Execute("301÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "695÷6=", 2)
